$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H4").Value = 219.70589
$ws_ALC.Range("I4").Value = 103.75
$ws_ALC.Range("J4").Value = 498
$ws_ALC.Range("K4").Value = 103.75
$ws_ALC.Range("L4").Value = 498
$ws_ALC.Range("M4").Value = 10.25
$ws_ALC.Range("N4").Value = -726

$ws_ALC.Range("H17").Value = 1378208.8
$ws_ALC.Range("J17").Value = 1443791.2
$ws_ALC.Range("L17").Value = 4331373.6
$ws_ALC.Range("N17").Value = -4331709.6

$ws_ALC.Range("H43").Value = 72917970
$ws_ALC.Range("I43").Value = 166667230
$ws_ALC.Range("J43").Value = 16668402
$ws_ALC.Range("K43").Value = 166667230
$ws_ALC.Range("L43").Value = 16668402
$ws_ALC.Range("M43").Value = -166667161
$ws_ALC.Range("N43").Value = -16668540

$ws_ALC.Range("H129").Value = 464965.1
$ws_ALC.Range("I129").Value = 565.7857
$ws_ALC.Range("J129").Value = 573324.94
$ws_ALC.Range("K129").Value = 1697.3571
$ws_ALC.Range("L129").Value = 1719974.82
$ws_ALC.Range("M129").Value = 3302.6429
$ws_ALC.Range("N129").Value = -1729974.82

$ws_ALC.Range("H132").Value = 2802954
$ws_ALC.Range("I132").Value = 3970362.2
$ws_ALC.Range("J132").Value = 1174.5333
$ws_ALC.Range("K132").Value = 11911086.6
$ws_ALC.Range("L132").Value = 3523.5999
$ws_ALC.Range("M132").Value = -11908556.6
$ws_ALC.Range("N132").Value = -8583.599900000001

$ws_ALC.Range("H135").Value = 654.05
$ws_ALC.Range("I135").Value = 469.17545
$ws_ALC.Range("K135").Value = 4222.57905
$ws_ALC.Range("M135").Value = -1687.57905

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 1169.9642
$ws_ARM.Range("I2").Value = 896.86664
$ws_ARM.Range("K2").Value = 896.86664
$ws_ARM.Range("M2").Value = -783.86664

$ws_ARM.Range("H3").Value = 100
$ws_ARM.Range("I3").Value = 100
$ws_ARM.Range("K3").Value = 100
$ws_ARM.Range("M3").Value = 15

$ws_ARM.Range("H61").Value = 3356.389
$ws_ARM.Range("I61").Value = 2938.875
$ws_ARM.Range("J61").Value = 3690.4
$ws_ARM.Range("K61").Value = 2938.875
$ws_ARM.Range("L61").Value = 3690.4
$ws_ARM.Range("M61").Value = -2726.875
$ws_ARM.Range("N61").Value = -4114.4

$ws_ARM.Range("H116").Value = 1169.9642
$ws_ARM.Range("I116").Value = 896.86664
$ws_ARM.Range("K116").Value = 896.86664
$ws_ARM.Range("M116").Value = 1397.13336

$ws_ARM.Range("H136").Value = 3356.389
$ws_ARM.Range("I136").Value = 2938.875
$ws_ARM.Range("J136").Value = 3690.4
$ws_ARM.Range("K136").Value = 8816.625
$ws_ARM.Range("L136").Value = 11071.2
$ws_ARM.Range("M136").Value = -6266.625
$ws_ARM.Range("N136").Value = -16171.2

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 1169.9642
$ws_BSM.Range("I3").Value = 896.86664
$ws_BSM.Range("K3").Value = 896.86664
$ws_BSM.Range("M3").Value = -782.86664

$ws_BSM.Range("H124").Value = 57500
$ws_BSM.Range("J124").Value = 57500
$ws_BSM.Range("L124").Value = 57500
$ws_BSM.Range("N124").Value = -67320

$ws_BSM.Range("H134").Value = 56243.05
$ws_BSM.Range("I134").Value = 172018.67
$ws_BSM.Range("J134").Value = 2808.1538
$ws_BSM.Range("K134").Value = 516056.01
$ws_BSM.Range("L134").Value = 8424.4614
$ws_BSM.Range("M134").Value = -513521.01
$ws_BSM.Range("N134").Value = -13494.4614

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 5211413.5
$ws_CRP.Range("I31").Value = 2548
$ws_CRP.Range("J31").Value = 12824370
$ws_CRP.Range("K31").Value = 2548
$ws_CRP.Range("L31").Value = 12824370
$ws_CRP.Range("M31").Value = -2253
$ws_CRP.Range("N31").Value = -12824960

$ws_CRP.Range("H34").Value = 5211413.5
$ws_CRP.Range("I34").Value = 2548
$ws_CRP.Range("J34").Value = 12824370
$ws_CRP.Range("K34").Value = 2548
$ws_CRP.Range("L34").Value = 12824370
$ws_CRP.Range("M34").Value = -2346
$ws_CRP.Range("N34").Value = -12824774

$ws_CRP.Range("H58").Value = 8547911
$ws_CRP.Range("I58").Value = 903.24243
$ws_CRP.Range("J58").Value = 55556456
$ws_CRP.Range("K58").Value = 903.24243
$ws_CRP.Range("L58").Value = 55556456
$ws_CRP.Range("M58").Value = -700.24243
$ws_CRP.Range("N58").Value = -55556862

$ws_CRP.Range("H99").Value = 3015.524
$ws_CRP.Range("I99").Value = 3034.6667
$ws_CRP.Range("J99").Value = 3001.1667
$ws_CRP.Range("K99").Value = 3034.6667
$ws_CRP.Range("L99").Value = 3001.1667
$ws_CRP.Range("M99").Value = -1536.6667
$ws_CRP.Range("N99").Value = -5997.1667

$ws_CRP.Range("H126").Value = 3015.524
$ws_CRP.Range("I126").Value = 3034.6667
$ws_CRP.Range("J126").Value = 3001.1667
$ws_CRP.Range("K126").Value = 9104.000100000001
$ws_CRP.Range("L126").Value = 9003.500100000001
$ws_CRP.Range("M126").Value = -6634.000100000001
$ws_CRP.Range("N126").Value = -13943.5001

$ws_CRP.Range("H136").Value = 8547911
$ws_CRP.Range("I136").Value = 903.24243
$ws_CRP.Range("J136").Value = 55556456
$ws_CRP.Range("K136").Value = 2709.72729
$ws_CRP.Range("L136").Value = 166669368
$ws_CRP.Range("M136").Value = -159.7272899999998
$ws_CRP.Range("N136").Value = -166674468

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H131").Value = 4509868
$ws_CUL.Range("I131").Value = 7635.7144
$ws_CUL.Range("J131").Value = 8217588.5
$ws_CUL.Range("K131").Value = 22907.1432
$ws_CUL.Range("L131").Value = 24652765.5
$ws_CUL.Range("M131").Value = -17867.1432
$ws_CUL.Range("N131").Value = -24662845.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H123").Value = 24411.143
$ws_GSM.Range("J123").Value = 24411.143
$ws_GSM.Range("L123").Value = 24411.143
$ws_GSM.Range("N123").Value = -29311.143

$ws_GSM.Range("H132").Value = 40289.5
$ws_GSM.Range("I132").Value = 55170.95
$ws_GSM.Range("J132").Value = 8873.111
$ws_GSM.Range("K132").Value = 165512.85
$ws_GSM.Range("L132").Value = 26619.333
$ws_GSM.Range("M132").Value = -162982.85
$ws_GSM.Range("N132").Value = -31679.333

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H68").Value = 2084.4443
$ws_LTW.Range("I68").Value = 1190
$ws_LTW.Range("J68").Value = 2800
$ws_LTW.Range("K68").Value = 1190
$ws_LTW.Range("L68").Value = 2800
$ws_LTW.Range("M68").Value = -441
$ws_LTW.Range("N68").Value = -4298

$ws_LTW.Range("H71").Value = 2084.4443
$ws_LTW.Range("I71").Value = 1190
$ws_LTW.Range("J71").Value = 2800
$ws_LTW.Range("K71").Value = 5950
$ws_LTW.Range("L71").Value = 14000
$ws_LTW.Range("M71").Value = -2206
$ws_LTW.Range("N71").Value = -21488

$ws_LTW.Range("H122").Value = 3142
$ws_LTW.Range("I122").Value = 4053.5
$ws_LTW.Range("J122").Value = 2686.25
$ws_LTW.Range("K122").Value = 12160.5
$ws_LTW.Range("L122").Value = 8058.75
$ws_LTW.Range("M122").Value = -9710.5
$ws_LTW.Range("N122").Value = -12958.75

$ws_LTW.Range("H136").Value = 3095.3403
$ws_LTW.Range("I136").Value = 3130.6191
$ws_LTW.Range("J136").Value = 2799
$ws_LTW.Range("K136").Value = 9391.8573
$ws_LTW.Range("L136").Value = 8397
$ws_LTW.Range("M136").Value = -6841.8573
$ws_LTW.Range("N136").Value = -13497

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 1594.1666
$ws_WVR.Range("I81").Value = 1113
$ws_WVR.Range("J81").Value = 4000
$ws_WVR.Range("K81").Value = 2226
$ws_WVR.Range("L81").Value = 8000
$ws_WVR.Range("M81").Value = -1165
$ws_WVR.Range("N81").Value = -10122

$ws_WVR.Range("H84").Value = 1594.1666
$ws_WVR.Range("I84").Value = 1113
$ws_WVR.Range("J84").Value = 4000
$ws_WVR.Range("K84").Value = 11130
$ws_WVR.Range("L84").Value = 40000
$ws_WVR.Range("M84").Value = -5826
$ws_WVR.Range("N84").Value = -50608

$ws_WVR.Range("H126").Value = 8532.533
$ws_WVR.Range("I126").Value = 8532.533
$ws_WVR.Range("J126").Value = 0
$ws_WVR.Range("K126").Value = 25597.599
$ws_WVR.Range("L126").Value = 0
$ws_WVR.Range("M126").Value = -23127.599
$ws_WVR.Range("N126").ClearContents()

$ws_WVR.Range("H132").Value = 1005.873
$ws_WVR.Range("I132").Value = 926.2222
$ws_WVR.Range("J132").Value = 1483.7778
$ws_WVR.Range("K132").Value = 2778.6666
$ws_WVR.Range("L132").Value = 4451.3334
$ws_WVR.Range("M132").Value = -248.6666
$ws_WVR.Range("N132").Value = -9511.3334

$ws_WVR.Range("H136").Value = 2441.4722
$ws_WVR.Range("I136").Value = 2978.0962
$ws_WVR.Range("J136").Value = 1046.25
$ws_WVR.Range("K136").Value = 8934.2886
$ws_WVR.Range("L136").Value = 3138.75
$ws_WVR.Range("M136").Value = -6384.2886
$ws_WVR.Range("N136").Value = -8238.75
